$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pre-Alert Template Import")

$ws.Range("A3").Value = "JSSO1000248"
$ws.Range("B3").Value = "JSSO1000248"
$ws.Range("C3").Value = "JSSO1000248"
$ws.Range("AJ3").Value = "JSCN1000248"
$ws.Range("AN3").Value = "MBLJSSO1000248"
$ws.Range("AO3").Value = "HBLJSSO1000248"
